$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("G3").Value = "num_screen/max(undx,num_screen)"
$ws.Range("G5").Value = "num_diag/max(scr,num_diag)"
$ws.Range("G7").Value = "num_initiate/max(dx,num_initiate)"
$ws.Range("G9").Value = "num_loss/max(tx,num_loss)"
$ws.Activate()
$ws.Range("G10").Select()
